# Updating the base model for Astro
# - Shift all interval timestamps (column A) and "Lookup" labels (column D)
#   forward by 6 days: 19.06.2024 -> 25.06.2024
# - Refresh the Prediction values (column C) for rows 27-86 with the new
#   model output

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @(
    @{Row=2; Date=45468},
    @{Row=3; Date=45468.01041666666},
    @{Row=4; Date=45468.02083333334},
    @{Row=5; Date=45468.03125},
    @{Row=6; Date=45468.04166666666},
    @{Row=7; Date=45468.05208333334},
    @{Row=8; Date=45468.0625},
    @{Row=9; Date=45468.07291666666},
    @{Row=10; Date=45468.08333333334},
    @{Row=11; Date=45468.09375},
    @{Row=12; Date=45468.10416666666},
    @{Row=13; Date=45468.11458333334},
    @{Row=14; Date=45468.13541666666},
    @{Row=15; Date=45468.14583333334},
    @{Row=16; Date=45468.15625},
    @{Row=17; Date=45468.16666666666},
    @{Row=18; Date=45468.17708333334},
    @{Row=19; Date=45468.1875},
    @{Row=20; Date=45468.19791666666},
    @{Row=21; Date=45468.20833333334},
    @{Row=22; Date=45468.21875},
    @{Row=23; Date=45468.22916666666},
    @{Row=24; Date=45468.23958333334},
    @{Row=25; Date=45468.25},
    @{Row=26; Date=45468.26041666666},
    @{Row=27; Date=45468.27083333334},
    @{Row=28; Date=45468.28125},
    @{Row=29; Date=45468.29166666666},
    @{Row=30; Date=45468.30208333334},
    @{Row=31; Date=45468.3125},
    @{Row=32; Date=45468.32291666666},
    @{Row=33; Date=45468.33333333334},
    @{Row=34; Date=45468.34375},
    @{Row=35; Date=45468.35416666666},
    @{Row=36; Date=45468.36458333334},
    @{Row=37; Date=45468.375},
    @{Row=38; Date=45468.38541666666},
    @{Row=39; Date=45468.39583333334},
    @{Row=40; Date=45468.40625},
    @{Row=41; Date=45468.41666666666},
    @{Row=42; Date=45468.42708333334},
    @{Row=43; Date=45468.4375},
    @{Row=44; Date=45468.44791666666},
    @{Row=45; Date=45468.45833333334},
    @{Row=46; Date=45468.46875},
    @{Row=47; Date=45468.47916666666},
    @{Row=48; Date=45468.48958333334},
    @{Row=49; Date=45468.5},
    @{Row=50; Date=45468.51041666666},
    @{Row=51; Date=45468.52083333334},
    @{Row=52; Date=45468.53125},
    @{Row=53; Date=45468.54166666666},
    @{Row=54; Date=45468.55208333334},
    @{Row=55; Date=45468.5625},
    @{Row=56; Date=45468.57291666666},
    @{Row=57; Date=45468.58333333334},
    @{Row=58; Date=45468.59375},
    @{Row=59; Date=45468.60416666666},
    @{Row=60; Date=45468.61458333334},
    @{Row=61; Date=45468.625},
    @{Row=62; Date=45468.63541666666},
    @{Row=63; Date=45468.64583333334},
    @{Row=64; Date=45468.65625},
    @{Row=65; Date=45468.66666666666},
    @{Row=66; Date=45468.67708333334},
    @{Row=67; Date=45468.6875},
    @{Row=68; Date=45468.69791666666},
    @{Row=69; Date=45468.70833333334},
    @{Row=70; Date=45468.71875},
    @{Row=71; Date=45468.72916666666},
    @{Row=72; Date=45468.73958333334},
    @{Row=73; Date=45468.75},
    @{Row=74; Date=45468.76041666666},
    @{Row=75; Date=45468.77083333334},
    @{Row=76; Date=45468.78125},
    @{Row=77; Date=45468.79166666666},
    @{Row=78; Date=45468.80208333334},
    @{Row=79; Date=45468.8125},
    @{Row=80; Date=45468.82291666666},
    @{Row=81; Date=45468.83333333334},
    @{Row=82; Date=45468.84375},
    @{Row=83; Date=45468.85416666666},
    @{Row=84; Date=45468.86458333334},
    @{Row=85; Date=45468.875},
    @{Row=86; Date=45468.88541666666},
    @{Row=87; Date=45468.89583333334},
    @{Row=88; Date=45468.90625},
    @{Row=89; Date=45468.91666666666},
    @{Row=90; Date=45468.92708333334},
    @{Row=91; Date=45468.9375},
    @{Row=92; Date=45468.94791666666},
    @{Row=93; Date=45468.95833333334},
    @{Row=94; Date=45468.96875},
    @{Row=95; Date=45468.97916666666},
    @{Row=96; Date=45468.98958333334}
)

foreach ($item in $newDates) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Date
}

for ($r = 2; $r -le 96; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $oldLabel = $cell.Value2
    if ($oldLabel -ne $null) {
        $newLabel = $oldLabel -replace "^19\.06\.2024", "25.06.2024"
        $cell.Value = $newLabel
    }
}

$newPredictions = @(
    @{Row=27; Val=0.011},
    @{Row=28; Val=0.023},
    @{Row=29; Val=0.041},
    @{Row=30; Val=0.055},
    @{Row=31; Val=0.07000000000000001},
    @{Row=32; Val=0.111},
    @{Row=33; Val=0.166},
    @{Row=34; Val=0.202},
    @{Row=35; Val=0.204},
    @{Row=36; Val=0.193},
    @{Row=37; Val=0.234},
    @{Row=38; Val=0.327},
    @{Row=39; Val=0.394},
    @{Row=40; Val=0.431},
    @{Row=41; Val=0.465},
    @{Row=42; Val=0.497},
    @{Row=43; Val=0.529},
    @{Row=44; Val=0.556},
    @{Row=45; Val=0.579},
    @{Row=46; Val=0.6},
    @{Row=47; Val=0.616},
    @{Row=48; Val=0.6840000000000001},
    @{Row=49; Val=0.703},
    @{Row=50; Val=0.717},
    @{Row=51; Val=0.732},
    @{Row=52; Val=0.733},
    @{Row=53; Val=0.733},
    @{Row=54; Val=0.733},
    @{Row=55; Val=0.73},
    @{Row=56; Val=0.726},
    @{Row=57; Val=0.715},
    @{Row=58; Val=0.7},
    @{Row=59; Val=0.673},
    @{Row=60; Val=0.646},
    @{Row=61; Val=0.598},
    @{Row=62; Val=0.551},
    @{Row=63; Val=0.506},
    @{Row=64; Val=0.491},
    @{Row=65; Val=0.456},
    @{Row=66; Val=0.407},
    @{Row=67; Val=0.367},
    @{Row=68; Val=0.326},
    @{Row=69; Val=0.32},
    @{Row=70; Val=0.303},
    @{Row=71; Val=0.285},
    @{Row=72; Val=0.273},
    @{Row=73; Val=0.241},
    @{Row=74; Val=0.197},
    @{Row=75; Val=0.161},
    @{Row=76; Val=0.127},
    @{Row=77; Val=0.11},
    @{Row=78; Val=0.08500000000000001},
    @{Row=79; Val=0.074},
    @{Row=80; Val=0.066},
    @{Row=81; Val=0.054},
    @{Row=82; Val=0.046},
    @{Row=83; Val=0.038},
    @{Row=84; Val=0.025},
    @{Row=85; Val=0.017},
    @{Row=86; Val=0.011}
)

foreach ($item in $newPredictions) {
    $ws.Cells.Item($item.Row, 3).Value = $item.Val
}
